$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cluster names to remove entirely (outbreaks no longer listed)
$removeNames = @(
    'Al Haj Halal Meats Glenroy',
    'Chemist Warehouse Campbellfield DC',
    'City of Wyndham Community',
    'Coles Campbellfield Plaza Campbellfield',
    'Construction Site 1 Warde Street Footscray',
    'Direct Freight Express Campbellfield',
    'Don Watson Coldstore Derrimut',
    'Epworth Healthcare Epworth Richmond Emergency Department',
    'Green Leaves Early Learning Cairnlea',
    'Kippers Seafood Werribee',
    'National Gallery of Victoria Melbourne',
    'Oscar Romero Catholic Primary School Craigieburn',
    'Sharpline Stainless Steel Coburg North',
    'Tek Foods Somerton',
    'The Huntly-Goornong Rail Works',
    'Yara Childcare Centre Truganina'
)

# Find the last used row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Collect row indices whose cluster name is in the removal list
$rowsToDelete = New-Object System.Collections.ArrayList
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($removeNames -contains $name) {
        [void]$rowsToDelete.Add($r)
    }
}

# Delete from bottom to top so row numbers of not-yet-deleted rows stay valid
for ($i = $rowsToDelete.Count - 1; $i -ge 0; $i--) {
    $r = $rowsToDelete[$i]
    $ws.Rows($r).Delete()
}

# Updated Active cases counts for the remaining clusters
$updates = @(
    @{Name='126 Racecourse Road Public Housing Tower Flemington'; Value=7},
    @{Name='139 Highett St Apartment Complex Richmond'; Value=9},
    @{Name='3175 The Bays Aged Care Facility Hastings'; Value=8},
    @{Name='3535 Opal Meadow Heights Aged Care Community Meadow Heights'; Value=11},
    @{Name='Al-Taqwa College Truganina'; Value=6},
    @{Name='Allbright Manor Aged Care Croydon North Tier 1A'; Value=5},
    @{Name='Australia Post Distribution Centre Sunshine West'; Value=5},
    @{Name='Australian Lamb Colac East'; Value=9},
    @{Name='Baker Bleu Caulfield North'; Value=7},
    @{Name='Baxter Foods Australia Campbellfield'; Value=5},
    @{Name='CFMEU Melbourne Office'; Value=5},
    @{Name='CS Square Caroline Springs'; Value=13},
    @{Name='Cafe Roco Dandenong'; Value=6},
    @{Name='Campbellfield Ford Complex Vaccination Clinic Campbellfield'; Value=10},
    @{Name='Cardinia Lakes Early Learning Centre Pakenham'; Value=5},
    @{Name='Caroline Springs Police Station'; Value=5},
    @{Name='Carton Finishing Pty. Ltd. Campbellfield'; Value=12},
    @{Name='Chemist Warehouse Fillo Drive Somerton'; Value=8},
    @{Name='Coles Coburg North Village'; Value=10},
    @{Name='Coles Pakenham Place Shopping Centre'; Value=5},
    @{Name='Coles Roxburgh Village Roxburgh Park'; Value=5},
    @{Name='Community Kids Bayswater Early Education Centre Bayswater North'; Value=18},
    @{Name='Construction Site Olea Apartment Caulfield North'; Value=7},
    @{Name='Costco Wholesale Epping'; Value=19},
    @{Name='Crusader Caravans Epping'; Value=17},
    @{Name='Dandenong Police Station Dandenong'; Value=9},
    @{Name='DayHab Rehabilitation Treatment Centre Ringwood East'; Value=7},
    @{Name='Disability Residence Life without Barriers Ashwood'; Value=5},
    @{Name='Ermha365 Residential Disability Care Services Paperbark St Doveton'; Value=9},
    @{Name='FedEx Station Melbourne Airport'; Value=14},
    @{Name='Fine Food Holdings Pty Ltd Dandenong South'; Value=10},
    @{Name='Fonterra Manufacturing Workplace Campbellfield'; Value=9},
    @{Name='General Foods Campbellfield'; Value=13},
    @{Name='Gladstone Parade Early Learning & Kinder Glenroy'; Value=7},
    @{Name='Goodstart Early Learning Altona'; Value=7},
    @{Name='Green Leaves Early Learning Centre Highlands Craigieburn'; Value=16},
    @{Name='Greenvale Primary School'; Value=5},
    @{Name='HEI Schools Emerald Early Learning Centre Emerald'; Value=5},
    @{Name='Hamilton Marino 236 Jasper Road McKinnon'; Value=7},
    @{Name='Hello Fresh Warehouse Ravenhall'; Value=7},
    @{Name='ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine'; Value=10},
    @{Name='Ibis Kingsgate Hotel Melbourne'; Value=5},
    @{Name='Industrial Galvanizers Valmont Coatings Campbellfield'; Value=8},
    @{Name='Inghams Enterprises Thomastown'; Value=6},
    @{Name='Kool Kidz Childcare Narre Warren'; Value=15},
    @{Name='Lantmannen Unibake Australia Mordialloc'; Value=7},
    @{Name='Linfox Somerton National Distribution Centre Somerton'; Value=6},
    @{Name='Mecca Distribution Centre Warehouse Melbourne Airport'; Value=7},
    @{Name='Melbourne Assessment Prison West Melbourne'; Value=9},
    @{Name='Melbourne Metropolitan Remand Centre Ravenhall'; Value=8},
    @{Name='Melbourne West Police Station Docklands'; Value=8},
    @{Name='Mill Park Police Station Mill Park'; Value=9},
    @{Name='MyCentre Childcare Broadmeadows'; Value=8},
    @{Name='Nido Early School Ascot Vale'; Value=26},
    @{Name='Nido Early School Glenroy'; Value=23},
    @{Name='Northern Health Northern Hospital Epping Emergency Department Tier 1B'; Value=39},
    @{Name='Northern Health The Northern Hospital Epping'; Value=20},
    @{Name='OnQ Plumbing and Excavations Craigieburn'; Value=8},
    @{Name='Oporto Coolaroo'; Value=5},
    @{Name='Our Lady Help of Christian''s Primary School Brunswick East'; Value=7},
    @{Name='Pacific Meat Thomastown'; Value=5},
    @{Name='Private Residence Daycare Allumba Way Wollert'; Value=8},
    @{Name='Ravenhall Correctional Centre Ravenhall'; Value=10},
    @{Name='Richmond Quarter 261-271 Bridge Road Construction Site Richmond'; Value=7},
    @{Name='Sacca''s Fruit World Broadmeadows Central Shopping Centre'; Value=6},
    @{Name='St Margaret''s Primary School OSHC Maribyrnong'; Value=12},
    @{Name='St Vincents Hospital Emergency Department Melbourne'; Value=18},
    @{Name='The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'; Value=15},
    @{Name='The Royal Melbourne Hospital Parkville'; Value=6},
    @{Name='The Royal Melbourne Hospital Parkville Emergency Department'; Value=7},
    @{Name='The Royal Melbourne Hospital Ward 6SE Parkville'; Value=10},
    @{Name='The Royal Talbot Rehabilitation Centre Kew'; Value=12},
    @{Name='ThorwestenCabinets Pakenham'; Value=10},
    @{Name='Truganina Early Learning Centre Truganina'; Value=8},
    @{Name='Visy Recycling Springvale'; Value=28},
    @{Name='Wallaby Childcare Wollert'; Value=18},
    @{Name='Werribee Mercy Hospital Emergency Department'; Value=14},
    @{Name='Western Health Footscray Hospital Emergency Department'; Value=7},
    @{Name='Western Health Sunshine Hospital Emergency Department'; Value=14},
    @{Name='Western Health Sunshine Hospital GEM Ward St Albans'; Value=6}
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    foreach ($u in $updates) {
        if ($u.Name -eq $name) {
            $ws.Cells.Item($r, 2).Value = $u.Value
            break
        }
    }
}
